# Updates cryptos list prices / 1h volume deltas (and re-ranks a few coins
# in rows 46-51, inserting BabyDogeCoin while the Price column keeps its
# text formatting). Cells whose new Price value looks like a plain decimal
# number (e.g. '214.18') are given a leading apostrophe so Excel stores
# them as text rather than re-typing them as a Number, matching the
# original column formatting (the thousands-style prices like
# '25.832.97' already stay text on their own since they aren't valid
# numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.832.97'
$ws.Range('D3').Value = '1.630.22'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').Value = '''214.18'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '''19.52'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').Value = '''0.0791'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '1.856.36'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').Value = '1.589.95'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '''0.543'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '0.0₃0752'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '''62.57'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '25.840.64'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = '''192.87'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').Value = '''143.30'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('D27').Value = '''0.125'
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('D28').Value = '''6.82'
$ws.Range('D29').Value = '''15.39'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('D36').Value = '''0.900'
$ws.Range('D37').Value = '1.136.74'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '''98.94'
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('D44').Value = '''0.792'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '1.766.12'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0113'
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '''56.18'
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0527'
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.44'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.415'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''7.62'
$ws.Range('E51').Value = '  +1.13%  '
